# Gallery_AttachedFunctionality_FIM_Node.xlsx
# Adds three new market test-data sheets (Netherlands, Austria, Denmark),
# modelled on the existing "Greece" sheet, and leaves the workbook with
# Austria as the active / selected sheet (matching the authored diff).

$wb = $excel.ActiveWorkbook
$greece = $wb.Worksheets.Item("Greece")

# ---------------------------------------------------------------------
# Netherlands  (plain copy of the Greece template)
# ---------------------------------------------------------------------
$greece.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$netherlands = $wb.Worksheets.Item($wb.Worksheets.Count)
$netherlands.Name = "Netherlands"
$netherlands.Range("B2").Value = "Netherlands Market"
$netherlands.Range("B4").Value = "NGC-3144/T2196"
$netherlands.Cells.Select() | Out-Null

# ---------------------------------------------------------------------
# Austria  (copy of the Greece template, with four extra "Attached
# Functionality" rows inserted before the trailing Wg/Attached rows)
# ---------------------------------------------------------------------
$greece.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$austria = $wb.Worksheets.Item($wb.Worksheets.Count)
$austria.Name = "Austria"
$austria.Range("B2").Value = "Austria Market"

$austria.Rows("10:13").Insert()
$austria.Range("A9").Copy()
$austria.Range("A10:A13").PasteSpecial(-4122)   # xlPasteFormats
$austria.Range("A10").Value = "Multichannel Transmission Unit"
$austria.Range("A11").Value = "Transmission Unit"
$austria.Range("A12").Value = "Transmission Unit and Keysafe"
$austria.Range("A13").Value = "Multichannel Transmission Unit"

$austria.Range("B4").Value = "NGC-3817/T2313"
$austria.Range("C11").Select() | Out-Null

# ---------------------------------------------------------------------
# Denmark  (plain copy of the Greece template)
# ---------------------------------------------------------------------
$greece.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-2913/T2196"
$denmark.Range("B4").Select() | Out-Null

# Austria ends up the active sheet, scrolled so it is visible in the tab
# strip (matches firstSheet="3"/activeTab="12" in the authored workbook).
$austria.Activate() | Out-Null
$austria.Range("C11").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(3) | Out-Null
